# Import-Format.xlsx edit script
# - Adds a "Content Type ID (contentTypeId)" column
# - Drops unused columns (Card Image, Intro Text, Orange Text, Source of information, News/Inspiration)
# - Refreshes the sample/demo data rows and adds a second demo row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove obsolete columns (rightmost first so letters keep their meaning)
#    E = Card Image, G = Intro Text, I = Orange Text, J = Source of information,
#    L = News/Inspiration
# ---------------------------------------------------------------------------
$ws.Columns("L").Delete()
$ws.Columns("J").Delete()
$ws.Columns("I").Delete()
$ws.Columns("G").Delete()
$ws.Columns("E").Delete()

# After the deletions the remaining columns are:
# A Title | B Slug | C Published Date | D Excerpt | E Featured Image | F Editor | G Tags

# ---------------------------------------------------------------------------
# 2. Add the new "Content Type ID" header in column H
# ---------------------------------------------------------------------------
$hHeader = $ws.Range("H1")
$hHeader.Value = "Content Type ID`n(contentTypeId)"
$hHeader.Font.Bold = $true
$hHeader.Interior.Color = $ws.Range("G1").Interior.Color
$hHeader.WrapText = $true
$hHeader.VerticalAlignment = -4108

# Make the trailing "(contentTypeId)" run bold + smaller, like the rest of the
# workbook's secondary-label convention
$sub = $hHeader.Characters(17, 15)
$sub.Font.Bold = $true
$sub.Font.Size = 9

$ws.Rows("1").RowHeight = 43.5

# ---------------------------------------------------------------------------
# 3. Refresh row 2 with new demo content and drop the inherited header fill
# ---------------------------------------------------------------------------
$ws.Range("A2:H2").Interior.ColorIndex = -4142

$ws.Range("A2").Value = "Post Title"
$ws.Range("B2").Value = "post-title"
$ws.Range("C2").Value = 45916.5105324074
$ws.Range("D2").Value = "Lorem Ipsum"
$ws.Range("E2").Value = "https://picsum.photos/536/354"
$ws.Range("F2").Value = "<p>Editor</p>"
$ws.Range("G2").Value = "test 11,test 22"
$ws.Range("H2").Value = "news"

$ws.Range("A2").Font.Bold = $true
$ws.Range("B2").Font.Bold = $true
$ws.Range("D2").Font.Bold = $true
$ws.Range("F2").Font.Bold = $true

# ---------------------------------------------------------------------------
# 4. Add row 3 as a second demo row, mirroring row 2's formatting
# ---------------------------------------------------------------------------
$ws.Range("A2:H2").Copy()
$ws.Range("A3").PasteSpecial(-4104)

$ws.Range("A3").Value = "Demo Title"
$ws.Range("B3").Value = "demo-title"
$ws.Range("C3").Value = 45916.5105324074
$ws.Range("D3").Value = "Lorem Ipsum"
$ws.Range("E3").Value = "https://picsum.photos/536/354"
$ws.Range("F3").Value = "<p>Editor</p>"
$ws.Range("G3").Value = "test 11,test 22"
$ws.Range("H3").Value = "news"

# ---------------------------------------------------------------------------
# 5. Hyperlinks for the picture-url cells
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E2"), "https://picsum.photos/536/354", "", "", "https://picsum.photos/536/354")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://picsum.photos/536/354", "", "", "https://picsum.photos/536/354")

# ---------------------------------------------------------------------------
# 6. View / selection bookkeeping
# ---------------------------------------------------------------------------
$ws.Range("H3").Select()

Write-Host "done"
